# Slide 3, "Content Placeholder 2" shape, paragraph 14 (the "Configuration"
# sub-bullet that reads "tsc init") needs to become "Tsc --init", keeping the
# original three-run split (run1 = command name w/ spellcheck-err flag,
# run2 = separator/space, run3 = "init" w/ spellcheck-err flag).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(14)
$start = $para.Start

# Run 1: "tsc" (3 chars) -> "Tsc"
$tr.Characters($start, 3).Text = "Tsc"

# Run 2: " " (1 char) -> " --"
$tr.Characters($start + 3, 1).Text = " --"
